$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.9083094555873925
$ws.Range("B3").Value = 0.913961038961039
$ws.Range("B4").Value = 0.929042904290429
$ws.Range("B5").Value = 0.9214402618657938
